# Update Excel file with latest predictions
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Home win"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Home win")

$ws.Cells.Item(2, 1).Value = "09-01-2025 16:00"
$ws.Cells.Item(2, 2).Value = "BRAZIL"
$ws.Cells.Item(2, 3).Value = "SÃO PAULO YOUTH CUP"
$ws.Cells.Item(2, 4).Value = "São José PA U20 - Retrô U20"
$ws.Cells.Item(2, 5).Value = 70
$ws.Cells.Item(2, 6).Value = 2.2

$ws.Cells.Item(3, 1).Value = "10-01-2025 15:00"
$ws.Cells.Item(3, 2).Value = "PORTUGAL"
$ws.Cells.Item(3, 3).Value = "LIGA REVELAÇÃO U23"
$ws.Cells.Item(3, 4).Value = "Academico Viseu U23 - Sporting Braga U23"
$ws.Cells.Item(3, 5).Value = 73.3
$ws.Cells.Item(3, 6).Value = 2.1

$ws.Cells.Item(4, 1).Value = "10-01-2025 19:30"
$ws.Cells.Item(4, 2).Value = "SPAIN"
$ws.Cells.Item(4, 3).Value = "SEGUNDA DIVISIÓN"
$ws.Cells.Item(4, 4).Value = "Granada CF - Burgos"
$ws.Cells.Item(4, 5).Value = 80
$ws.Cells.Item(4, 6).Value = 1.73

# ---------------------------------------------------------------
# Sheet "Draw"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Draw")

$ws.Cells.Item(3, 1).Value = "10-01-2025 20:00"
$ws.Cells.Item(3, 2).Value = "FRANCE"
$ws.Cells.Item(3, 3).Value = "LIGUE 1"
$ws.Cells.Item(3, 4).Value = "Auxerre - Lille"
$ws.Cells.Item(3, 5).Value = 66.7
$ws.Cells.Item(3, 6).Value = 3.7

$ws.Cells.Item(4, 1).Value = "10-01-2025 17:00"
$ws.Cells.Item(4, 2).Value = "EGYPT"
$ws.Cells.Item(4, 3).Value = "PREMIER LEAGUE"
$ws.Cells.Item(4, 4).Value = "El Geish - Ceramica Cleopatra"
$ws.Cells.Item(4, 5).Value = 60
$ws.Cells.Item(4, 6).Value = 3

$ws.Cells.Item(5, 1).Value = "10-01-2025 18:30"
$ws.Cells.Item(5, 2).Value = "FRANCE"
$ws.Cells.Item(5, 3).Value = "NATIONAL 1"
$ws.Cells.Item(5, 4).Value = "Dijon - Chateauroux"
$ws.Cells.Item(5, 5).Value = 60
$ws.Cells.Item(5, 6).Value = 3.75

# ---------------------------------------------------------------
# Sheet "Btts"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Btts")

$ws.Cells.Item(3, 1).Value = "09-01-2025 18:15"
$ws.Cells.Item(3, 2).Value = "BRAZIL"
$ws.Cells.Item(3, 3).Value = "SÃO PAULO YOUTH CUP"
$ws.Cells.Item(3, 4).Value = "Tupã U20 - CRB U20"
$ws.Cells.Item(3, 5).Value = 88
$ws.Cells.Item(3, 6).Value = 1.8

$ws.Cells.Item(4, 1).Value = "10-01-2025 19:00"
$ws.Cells.Item(4, 2).Value = "NETHERLANDS"
$ws.Cells.Item(4, 3).Value = "EERSTE DIVISIE"
$ws.Cells.Item(4, 4).Value = "Den Bosch - Jong Utrecht"
$ws.Cells.Item(4, 5).Value = 86.7
$ws.Cells.Item(4, 6).Value = 1.7

$ws.Cells.Item(5, 1).Value = "10-01-2025 18:30"
$ws.Cells.Item(5, 2).Value = "FRANCE"
$ws.Cells.Item(5, 3).Value = "NATIONAL 1"
$ws.Cells.Item(5, 4).Value = "Valenciennes - Gobelins"
$ws.Cells.Item(5, 5).Value = 76
$ws.Cells.Item(5, 6).Value = 1.85

$ws.Cells.Item(6, 1).Value = "10-01-2025 19:30"
$ws.Cells.Item(6, 2).Value = "ITALY"
$ws.Cells.Item(6, 3).Value = "SERIE C - GIRONE B"
$ws.Cells.Item(6, 4).Value = "Lucchese - Spal"
$ws.Cells.Item(6, 5).Value = 86.7
$ws.Cells.Item(6, 6).Value = 1.83

# ---------------------------------------------------------------
# Sheet "Over_Under"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Over_Under")

$ws.Cells.Item(2, 1).Value = "09-01-2025 18:15"
$ws.Cells.Item(2, 4).Value = "Tupã U20 - CRB U20"
$ws.Cells.Item(2, 6).Value = 1.8
$ws.Cells.Item(2, 7).Value = 40
$ws.Cells.Item(2, 8).Value = 2.95

$ws.Cells.Item(4, 1).Value = "10-01-2025 19:45"
$ws.Cells.Item(4, 2).Value = "SCOTLAND"
$ws.Cells.Item(4, 3).Value = "CHAMPIONSHIP"
$ws.Cells.Item(4, 4).Value = "Hamilton Academical - Partick"
$ws.Cells.Item(4, 5).Value = 80
$ws.Cells.Item(4, 6).Value = 1.73
$ws.Cells.Item(4, 7).Value = 50
$ws.Cells.Item(4, 8).Value = 2.88

$ws.Cells.Item(5, 1).Value = "10-01-2025 14:00"
$ws.Cells.Item(5, 2).Value = "WORLD"
$ws.Cells.Item(5, 3).Value = "FRIENDLIES CLUBS"
$ws.Cells.Item(5, 4).Value = "Karlsruher SC - Servette FC"
$ws.Cells.Item(5, 5).Value = 80
$ws.Cells.Item(5, 6).Value = 1.5
$ws.Cells.Item(5, 7).Value = 73.3
$ws.Cells.Item(5, 8).Value = 2.25

# ---------------------------------------------------------------
# Sheet "Away Win"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Away Win")

$ws.Cells.Item(2, 1).Value = "09-01-2025 16:00"
$ws.Cells.Item(2, 2).Value = "BRAZIL"
$ws.Cells.Item(2, 3).Value = "SÃO PAULO YOUTH CUP"
$ws.Cells.Item(2, 4).Value = "Botafogo SP U20 - Tuna Luso PA U20"
$ws.Cells.Item(2, 5).Value = 70
$ws.Cells.Item(2, 6).Value = 10

$ws.Cells.Item(3, 1).Value = "09-01-2025 18:15"
$ws.Cells.Item(3, 2).Value = "BRAZIL"
$ws.Cells.Item(3, 3).Value = "SÃO PAULO YOUTH CUP"
$ws.Cells.Item(3, 4).Value = "Comercial De Tietê U20 - América RN U20"
$ws.Cells.Item(3, 5).Value = 70
$ws.Cells.Item(3, 6).Value = 2.85

$ws.Cells.Item(4, 1).Value = "10-01-2025 19:00"
$ws.Cells.Item(4, 2).Value = "FRANCE"
$ws.Cells.Item(4, 3).Value = "LIGUE 2"
$ws.Cells.Item(4, 4).Value = "Martigues - Clermont Foot"
$ws.Cells.Item(4, 5).Value = 80
$ws.Cells.Item(4, 6).Value = 1.75
